$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new columns before column D (shifts D:K -> F:M)
$ws.Range("D:E").Insert()

# Copy number formatting from column F (old column D, now shifted) onto the
# two newly inserted columns D:E so they pick up the same formats as the
# rest of the quarterly data (dates in row style, numbers in data style)
$ws.Range("F7:F102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$allRowData = @{
  7 = @(43465, 43373, 43281, 43190, 43100, 43008, 42916, 42825, 42735, 42643)
  8 = @(107200, 100500, 96600, 113400, 86000, 81500, 76100, 90300, 68700, 65000)
  9 = @(38100, 35900, 38700, 38700, 36800, 35000, 33200, 32100, 30400, 28400)
  10 = @(69100, 64600, 57900, 74700, 49200, 46500, 42900, 58200, 38300, 36600)
  11 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
  12 = @(12800, 11400, 10400, 9100, 9300, 8900, 7600, 7000, 7200, 7300)
  13 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  14 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  15 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  16 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
  17 = @(100200, 96700, 101800, 92900, 85900, 81000, 79500, 75400, 70300, 67500)
  18 = @(7000, 3800, -5200, 20500, 100, 500, -3400, 14900, -1600, -2500)
  19 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
  20 = @(400, 300, 400, 200, 100, 100, 0, -100, 0, 100)
  21 = @(15900, 12300, 4700, 27900, 7000, 7300, 3000, 20400, 3200, 1800)
  22 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  23 = @(7400, 4000, -4800, 20700, 300, 600, -3400, 14800, -1600, -2500)
  24 = @(1700, -5800, -3300, -18500, -200, 100, 500, 0, 0, 100)
  25 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  26 = @(5700, 9900, -1600, 39200, 400, 500, -3800, 14800, -1700, -2600)
  27 = @(5700, 9900, -1600, 39200, 400, 500, -3800, 14800, -1700, -2600)
  28 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  29 = @("NA", "NA", 0, 0, "NA", "NA", "NA", "NA", "NA", "NA")
  30 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  31 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  32 = @(-400, -300, -400, -200, -100, -100, 0, 100, 0, -100)
  33 = @(5700, 9900, -1600, 39200, 400, 500, -3800, 14800, -1700, -2600)
  34 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  35 = @(5700, 9900, -1600, 39200, 400, 500, -3800, 14800, -1700, -2600)
  38 = @(43465, 43373, 43281, 43190, 43100, 43008, 42916, 42825, 42735, 42643)
  39 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
  40 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
  41 = @(84100, 63700, 137200, 129500, 111000, 97400, 103500, 101500, 82300, 78000)
  42 = @(19900, 21200, "NA", "NA", "NA", "NA", "NA", 0, 0, 0)
  43 = @(21900, 19200, 3500, 3400, 2700, 2600, 2000, 2200, 2100, 1800)
  44 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  45 = @(1271300, 1182400, 1237600, 1364400, 1353200, 955500, 957300, 1183400, 1103800, 777400)
  46 = @(1397300, 1286500, 1378200, 1497400, 1466900, 1055400, 1062800, 1287100, 1188100, 857100)
  47 = @(64600, 59100, "NA", "NA", "NA", "NA", "NA", 0, 0, 0)
  48 = @(62700, 61800, 62000, 50400, 48400, 45000, 40800, 34300, 33600, 28700)
  49 = @(44600, 44100, 43700, 43200, 33000, 32500, 32300, 32000, 30600, 29100)
  50 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  51 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  52 = @(11200, 12700, 23600, 19900, 1100, 1300, 1500, 1600, 600, 500)
  53 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  54 = @(1580400, 1464100, 1507600, 1610900, 1549300, 1134300, 1137400, 1354900, 1252900, 915400)
  55 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
  56 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
  57 = @(3500, 3600, 3000, 2400, 1800, 1900, 2000, 1600, 3100, 2800)
  58 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  59 = @(1302000, 1203800, 1267900, 1383000, 1374800, 968700, 972800, 1197500, 1114600, 789200)
  60 = @(1305500, 1207400, 1270800, 1385400, 1376600, 970600, 974800, 1199100, 1117700, 792000)
  61 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  62 = @(23900, 24100, 23900, 21000, 14600, 15000, 15000, 10000, 9900, 4800)
  63 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  64 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  65 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  66 = @(1329400, 1231500, 1294800, 1406300, 1391200, 985500, 989800, 1209100, 1127600, 796800)
  67 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
  68 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  69 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  70 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  71 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  72 = @(61600, 55800, -6700, -5100, -44300, -44700, -45300, -41400, -56200, -54600)
  73 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  74 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  75 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  76 = @(250900, 232600, 212800, 204500, 158200, 148700, 147600, 145800, 125300, 118600)
  77 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  80 = @(43465, 43373, 43281, 43190, 43100, 43008, 42916, 42825, 42735, 42643)
  81 = @(5700, 9900, -1600, 39200, 400, 500, -3800, 14800, -1700, -2600)
  82 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
  83 = @(8600, 8200, 9600, 7200, 6800, 6700, 6300, 5600, 4800, 4300)
  84 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  85 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  86 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  87 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  88 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  89 = @(27000, 7300, 28500, 35200, 26000, 8200, 18600, 27900, 13500, 1900)
  90 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
  91 = @(-5100, -2400, -12300, -1400, -5300, -2700, -7400, -7700, -10500, -5800)
  92 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  93 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  94 = @(-98600, 29100, 100300, -15400, -416900, -7300, 214000, -85500, -332000, 462800)
  95 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
  96 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  97 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  98 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  99 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  100 = @(92100, -109900, -121100, -1200, 404500, -6900, -230700, 76800, 322800, -473200)
  101 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  102 = @(20500, -73500, 7700, 18500, 13600, -6100, 2000, 19200, 4300, -8500)
}

foreach ($r in $allRowData.Keys) {
    $vals = $allRowData[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = 4 + $i   # column D = 4
        $v = $vals[$i]
        if ($null -ne $v) {
            $ws.Cells.Item($r, $col).Value = $v
        }
    }
}
